# Update pl_clk_1 related divider values from 6 to 8 (187.5MHz)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C28").Value = 8
$ws.Range("C30").Value = 8

# Update selection on Sheet1 to reflect the last active cell (C28)
$ws.Activate()
$ws.Range("C28").Select()
